$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (borders etc.) from the last existing row (row 5) to the new row (row 6)
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Ensure the new cells are stored as text (matching the rest of the sheet),
# even though some values look numeric.
$ws.Range("A6:F6").NumberFormat = "@"

# Fill in the new row of data
$ws.Range("A6").Value = "20141364"
$ws.Range("B6").Value = "SARIWANGI TWIN 2X30S"
$ws.Range("C6").Value = "FES05N"
$ws.Range("D6").Value = "1"
$ws.Range("E6").Value = "7"
$ws.Range("F6").Value = "PT,(E-1B)"
